$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel;
# force them to remain Text so the literal formatted string is preserved.
$ws.Range("D2").Value = "64.750.13"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").Value = "3.169.74"
$ws.Range("E3").Value = "  -7.89%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.10"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.18"
$ws.Range("E6").Value = "  -5.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "3.166.57"
$ws.Range("E9").Value = "  -7.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.58"
$ws.Range("E11").Value = "  -5.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  -5.64%  "
$ws.Range("D13").Value = "3.721.06"
$ws.Range("E13").Value = "  -7.83%  "
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.40"
$ws.Range("D16").Value = "64.634.06"
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000164"
$ws.Range("E17").Value = "  -6.37%  "
$ws.Range("D18").Value = "3.169.55"
$ws.Range("E18").Value = "  -7.64%  "
$ws.Range("E19").Value = "  -4.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.99"
$ws.Range("E20").Value = "  -6.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "355.79"
$ws.Range("E21").Value = "  -4.75%  "
$ws.Range("E22").Value = "  -5.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.66"
$ws.Range("E24").Value = "  -6.42%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000118"
$ws.Range("E25").Value = "  -8.92%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.501"
$ws.Range("E26").Value = "  -7.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.90"
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.44"
$ws.Range("E32").Value = "  -7.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.03"
$ws.Range("E33").Value = "  -7.07%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.21"
$ws.Range("E34").Value = "  -6.13%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.68"
$ws.Range("E35").Value = "  -6.26%  "
$ws.Range("E36").Value = "  -8.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.42"
$ws.Range("E37").Value = "  -5.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.831"
$ws.Range("E38").Value = "  -5.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.34"
$ws.Range("E39").Value = "  -5.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.74"
$ws.Range("E40").Value = "  -4.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.53"
$ws.Range("E41").Value = "  -5.11%  "
$ws.Range("D42").Value = "2.659.24"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("E43").Value = "  -7.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.04"
$ws.Range("E44").Value = "  -5.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.16"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.18"
$ws.Range("E46").Value = "  -6.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0655"
$ws.Range("E47").Value = "  -6.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "323.34"
$ws.Range("E48").Value = "  -4.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0274"
$ws.Range("E49").Value = "  -4.93%  "
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("E51").Value = "  -0.07%  "
